# Atualizacao de bases das ligas, do dia: 17-06-2024 as 21:10
#
# The source data rows got re-sorted/re-matched against the result/odds
# records, which (for a handful of rows) swapped the full match record
# between two adjacent rows while the running "id" (column A), the
# league (column C) and the match Date (column D) stay attached to the
# original row number.
#
# Concretely: for each pair of rows below, every other column
# (B, E..AD : match id, teams, scores, odds, P/L, ...) is exchanged
# between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get swapped: B (match id) and E..AD (everything after the
# Date column). Columns A, C and D are left untouched.
$colsToSwap = @(2) + @(5..30)

$rowPairs = @(
    @(130, 131),
    @(143, 145),
    @(236, 237)
)

foreach ($rowPair in $rowPairs) {
    $rowA = $rowPair[0]
    $rowB = $rowPair[1]

    foreach ($col in $colsToSwap) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valueA = $cellA.Value2
        $valueB = $cellB.Value2

        # Only touch a cell when the incoming value actually differs from what
        # is already stored there, so untouched cells keep their original
        # (unmodified) representation.
        if ($valueA -ne $valueB) {
            $cellA.Value = $valueB
            $cellB.Value = $valueA
        }
    }
}
